$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The acquisition dates for rows 4 and 7 ("01/08/2022" and "05/04/2023")
# look like ambiguous dd/mm/yyyy dates that Excel's automatic type
# detection would otherwise convert into date serial numbers. Prefix
# them with an apostrophe so they are stored as literal text; the
# formatting copy further below (which reproduces data row 3's look)
# overwrites the resulting "quote prefix" marker so every cell in the
# new rows ends up sharing the exact same style.
$ws.Range("D4").Value = "'01/08/2022"
$ws.Range("D7").Value = "'05/04/2023"

# Clone the style (borders/font/alignment) and row height of the
# existing data row (row 3) onto the four new data rows (4-7).
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H7").PasteSpecial(-4122)
for ($r = 4; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = $ws.Rows.Item(3).RowHeight
}
$excel.CutCopyMode = 0

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 100111
$ws.Range("C4").Value = "Quadro branco magnético 2.00m x 1.20m"
$ws.Range("E4").Value = "NF-e 49123"
$ws.Range("F4").Value = "Coordenação do Curso de Pedagogia"
$ws.Range("G4").Value = "Irrecuperável"
$ws.Range("H4").Value = "Alienação/Leilão"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 100112
$ws.Range("C5").Value = "Cadeira universitária com prancheta fixa"
$ws.Range("D5").Value = "15/02/2020"
$ws.Range("E5").Value = "NF-e 18990"
$ws.Range("F5").Value = "Coordenação do Curso de História"
$ws.Range("G5").Value = "Irrecuperável"
$ws.Range("H5").Value = "Alienação/Leilão"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 100113
$ws.Range("C6").Value = "Longarina de 3 lugares para recepção, estofado azul"
$ws.Range("D6").Value = "30/05/2019"
$ws.Range("E6").Value = "NF-e 14321"
$ws.Range("F6").Value = "Vice-Reitoria"
$ws.Range("G6").Value = "Irrecuperável"
$ws.Range("H6").Value = "Alienação/Leilão"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 100114
$ws.Range("C7").Value = 'Televisor Smart 50" 4K LG'
$ws.Range("E7").Value = "NF-e 53112"
$ws.Range("F7").Value = "Pró-Reitoria de Graduação (PROGRAD)"
$ws.Range("G7").Value = "Irrecuperável"
$ws.Range("H7").Value = "Alienação/Leilão"
